$d = $word.ActiveDocument

# Title heading (appears twice: at top as Heading1, and later as bold text)
# Replace=2 (wdReplaceAll) replaces both occurrences in this single call.
$d.Content.Find.Execute(
    "Play Legend of Cleopatra Free - Review of Egyptian-Themed Slot",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Legend of Cleopatra for Free", 2)

# "What we like" list items
$d.Content.Find.Execute(
    "Well-done graphics and unique Egyptian theme",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Unique and unusual 6-reel layout", 2)

$d.Content.Find.Execute(
    "Unusual game engine with 6 reels and up to 100 paylines",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Up to 100 paylines for more chances to win", 2)

$d.Content.Find.Execute(
    "Golden wild symbol, double wild reels feature, and free spins feature",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Well-done graphics and a fitting soundtrack", 2)

$d.Content.Find.Execute(
    "Betting options for all devices from €0.20 to €100",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exciting extra features like Double Wild Reel and Free Spins", 2)

# "What we don't like" list items
$d.Content.Find.Execute(
    "May not appeal to players who are not interested in the ancient Egypt theme",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited betting options compared to some other slot games", 2)

$d.Content.Find.Execute(
    "Lacks bonus games beyond the free spins feature",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited selection of other titles in the same genre", 2)

# Italic summary text
$d.Content.Find.Execute(
    "Experience the mysteries of the pyramids and win big with Legend of Cleopatra. Play for free and read our review of this unique online slot game.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Legend of Cleopatra and play for free at your favorite online casino.", 2)
